$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "65.582.53"
$ws.Range("E2").Formula = "  -0.24%  "

# Row 3
$ws.Range("D3").Formula = "3.294.61"
$ws.Range("E3").Formula = "  +0.75%  "

# Row 4
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Formula = "  +0.13%  "

# Row 5
$ws.Range("D5").Formula = "'553.47"
$ws.Range("E5").Formula = "  -0.78%  "

# Row 6
$ws.Range("D6").Formula = "'182.45"
$ws.Range("E6").Formula = "  -1.22%  "

# Row 7
$ws.Range("D7").Formula = "'0.999"
$ws.Range("E7").Formula = "  -0.15%  "

# Row 8
$ws.Range("D8").Formula = "3.284.63"
$ws.Range("E8").Formula = "  +0.71%  "

# Row 9
$ws.Range("D9").Formula = "'0.569"
$ws.Range("E9").Formula = "  -3.37%  "

# Row 10
$ws.Range("D10").Formula = "'0.171"
$ws.Range("E10").Formula = "  -7.12%  "

# Row 11
$ws.Range("D11").Formula = "'0.570"
$ws.Range("E11").Formula = "  -2.22%  "

# Row 12
$ws.Range("D12").Formula = "'45.05"
$ws.Range("E12").Formula = "  -4.52%  "

# Row 13
$ws.Range("D13").Formula = "'0.0000258"
$ws.Range("E13").Formula = "  -3.15%  "

# Row 14
$ws.Range("D14").Formula = "3.832.78"
$ws.Range("E14").Formula = "  +0.94%  "

# Row 15
$ws.Range("D15").Formula = "'8.32"
$ws.Range("E15").Formula = "  -3.29%  "

# Row 16
$ws.Range("D16").Formula = "'562.93"
$ws.Range("E16").Formula = "  -11.12%  "

# Row 17
$ws.Range("D17").Formula = "65.529.32"
$ws.Range("E17").Formula = "  -0.26%  "

# Row 18
$ws.Range("E18").Formula = "  +0.24%  "

# Row 19
$ws.Range("D19").Formula = "3.299.83"
$ws.Range("E19").Formula = "  +0.91%  "

# Row 20
$ws.Range("D20").Formula = "'17.43"
$ws.Range("E20").Formula = "  -3.22%  "

# Row 21
$ws.Range("D21").Formula = "'10.69"
$ws.Range("E21").Formula = "  -5.66%  "

# Row 22
$ws.Range("D22").Formula = "'0.878"
$ws.Range("E22").Formula = "  -2.52%  "

# Row 23
$ws.Range("D23").Formula = "'17.55"
$ws.Range("E23").Formula = "  -3.99%  "

# Row 24
$ws.Range("D24").Formula = "'4.91"
$ws.Range("E24").Formula = "  +0.23%  "

# Row 25
$ws.Range("D25").Formula = "'96.97"
$ws.Range("E25").Formula = "  -8.95%  "

# Row 26
$ws.Range("D26").Formula = "'3.90"
$ws.Range("E26").Formula = "  -1.68%  "

# Row 27
$ws.Range("D27").Formula = "'5.94"
$ws.Range("E27").Formula = "  -0.13%  "

# Row 28
$ws.Range("D28").Formula = "'2.65"
$ws.Range("E28").Formula = "  -0.85%  "

# Row 29
$ws.Range("D29").Formula = "'9.18"
$ws.Range("E29").Formula = "  -3.35%  "

# Row 30
$ws.Range("D30").Formula = "'8.32"
$ws.Range("E30").Formula = "  -4.11%  "

# Row 31
$ws.Range("D31").Formula = "'30.06"
$ws.Range("E31").Formula = "  -0.54%  "

# Row 32
$ws.Range("D32").Formula = "'6.49"
$ws.Range("E32").Formula = "  +4.00%  "

# Row 33
$ws.Range("B33").Formula = "Bittensor"
$ws.Range("C33").Formula = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Formula = "'555.40"
$ws.Range("E33").Formula = "  +6.15%  "

# Row 34
$ws.Range("B34").Formula = "dogwifhat"
$ws.Range("C34").Formula = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Formula = "'3.61"
$ws.Range("E34").Formula = "  -8.47%  "

# Row 36
$ws.Range("D36").Formula = "3.751.48"
$ws.Range("E36").Formula = "  +0.79%  "

# Row 37
$ws.Range("B37").Formula = "Dai"
$ws.Range("C37").Formula = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Formula = "'0.999"
$ws.Range("E37").Formula = "  -0.10%  "

# Row 38
$ws.Range("B38").Formula = "Hedera"
$ws.Range("C38").Formula = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Formula = "'0.101"
$ws.Range("E38").Formula = "  -3.30%  "

# Row 39
$ws.Range("D39").Formula = "'55.44"
$ws.Range("E39").Formula = "  -3.70%  "

# Row 40
$ws.Range("D40").Formula = "'32.96"
$ws.Range("E40").Formula = "  +0.34%  "

# Row 41
$ws.Range("E41").Formula = "  -3.90%  "

# Row 42
$ws.Range("D42").Formula = "'3.09"
$ws.Range("E42").Formula = "  -9.15%  "

# Row 43
$ws.Range("B43").Formula = "ApeXProtocol"
$ws.Range("C43").Formula = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Formula = "'3.32"
$ws.Range("E43").Formula = "  +2.01%  "

# Row 44
$ws.Range("B44").Formula = "PEPE"
$ws.Range("C44").Formula = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Formula = "0.0₃0664"
$ws.Range("E44").Formula = "  -9.32%  "

# Row 45
$ws.Range("D45").Formula = "'2.52"
$ws.Range("E45").Formula = "  -6.86%  "

# Row 46
$ws.Range("D46").Formula = "'0.326"
$ws.Range("E46").Formula = "  -3.31%  "

# Row 47
$ws.Range("D47").Formula = "'0.0402"
$ws.Range("E47").Formula = "  -2.90%  "

# Row 48
$ws.Range("D48").Formula = "'2.96"
$ws.Range("E48").Formula = "  -12.76%  "

# Row 49
$ws.Range("D49").Formula = "'1.00"
$ws.Range("E49").Formula = "  +0.26%  "

# Row 50
$ws.Range("E50").Formula = "  -3.44%  "

# Row 51
$ws.Range("D51").Formula = "'2.47"
$ws.Range("E51").Formula = "  -4.89%  "
